$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (rows 2-4) ---
$ws.Range("B2").Value = 54.6
$ws.Range("C2").Value = 24.165

$ws.Range("B3").Value = 56.8
$ws.Range("C3").Value = 23.607

$ws.Range("B4").Value = 52.7
$ws.Range("C4").Value = 24.649

# --- Add new row 5 ---
# Build the text label "3" in a scratch cell far away so it is stored as
# text (not auto-converted to a number), then bring just the value over to
# A5 and overlay A4's formatting (bold/border/center) so A5 matches the
# style used by the other label cells in column A.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "3"
$scratch.Copy()
$ws.Range("A5").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$scratch.Clear() | Out-Null

$ws.Range("B5").Value = 56.2
$ws.Range("C5").Value = 23.398
